$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.405597999999999
$ws.Range("H2").Value = 25.216794
$ws.Range("I2").Value = 0.3214983278049074
$ws.Range("J2").Value = 0.3214983278049074
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 35.337883
$ws.Range("N2").Value = 106.013649
$ws.Range("O2").Value = 0.3968231145247413
$ws.Range("P2").Value = 0.3968231145247413
$ws.Range("Q2").Value = 297.036038669034
$ws.Range("R2").Value = 2673.324348021306
$ws.Range("S2").Value = 0.1275779677540396
$ws.Range("T2").Value = 0.1275779677540396

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.405597999999999
$ws.Range("H3").Value = 25.216794
$ws.Range("I3").Value = 0.3214983278049074
$ws.Range("J3").Value = 0.3214983278049074
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.49537033333333
$ws.Range("N3").Value = 121.486111
$ws.Range("O3").Value = 0.4547385869013757
$ws.Range("P3").Value = 0.4547385869013756
$ws.Range("Q3").Value = 340.387803883126
$ws.Range("R3").Value = 3063.490234948134
$ws.Range("S3").Value = 0.1461976952771588
$ws.Range("T3").Value = 0.1461976952771588

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.405597999999999
$ws.Range("H4").Value = 25.216794
$ws.Range("I4").Value = 0.3214983278049074
$ws.Range("J4").Value = 0.3214983278049074
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.218724
$ws.Range("N4").Value = 39.656172
$ws.Range("O4").Value = 0.1484382985738831
$ws.Range("P4").Value = 0.148438298573883
$ws.Range("Q4").Value = 111.111280016952
$ws.Range("R4").Value = 1000.001520152568
$ws.Range("S4").Value = 0.04772266477370898
$ws.Range("T4").Value = 0.04772266477370898

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.54617866666667
$ws.Range("H5").Value = 37.638536
$ws.Range("I5").Value = 0.4798677573772784
$ws.Range("J5").Value = 0.4798677573772784
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 35.337883
$ws.Range("N5").Value = 106.013649
$ws.Range("O5").Value = 0.3968231145247413
$ws.Range("P5").Value = 0.3968231145247413
$ws.Range("Q5").Value = 443.3553938197627
$ws.Range("R5").Value = 3990.198544377864
$ws.Range("S5").Value = 0.1904226180424546
$ws.Range("T5").Value = 0.1904226180424545

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.54617866666667
$ws.Range("H6").Value = 37.638536
$ws.Range("I6").Value = 0.4798677573772784
$ws.Range("J6").Value = 0.4798677573772784
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.49537033333333
$ws.Range("N6").Value = 121.486111
$ws.Range("O6").Value = 0.4547385869013757
$ws.Range("P6").Value = 0.4547385869013756
$ws.Range("Q6").Value = 508.0621513748329
$ws.Range("R6").Value = 4572.559362373496
$ws.Range("S6").Value = 0.2182143858892758
$ws.Range("T6").Value = 0.2182143858892757

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.54617866666667
$ws.Range("H7").Value = 37.638536
$ws.Range("I7").Value = 0.4798677573772784
$ws.Range("J7").Value = 0.4798677573772784
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.218724
$ws.Range("N7").Value = 39.656172
$ws.Range("O7").Value = 0.1484382985738831
$ws.Range("P7").Value = 0.148438298573883
$ws.Range("Q7").Value = 165.8444730493547
$ws.Range("R7").Value = 1492.600257444192
$ws.Range("S7").Value = 0.07123075344554813
$ws.Range("T7").Value = 0.07123075344554812

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.193298666666666
$ws.Range("H8").Value = 15.579896
$ws.Range("I8").Value = 0.1986339148178141
$ws.Range("J8").Value = 0.1986339148178141
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 35.337883
$ws.Range("N8").Value = 106.013649
$ws.Range("O8").Value = 0.3968231145247413
$ws.Range("P8").Value = 0.3968231145247413
$ws.Range("Q8").Value = 183.5201806667226
$ws.Range("R8").Value = 1651.681626000504
$ws.Range("S8").Value = 0.07882252872824716
$ws.Range("T8").Value = 0.07882252872824716

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.193298666666666
$ws.Range("H9").Value = 15.579896
$ws.Range("I9").Value = 0.1986339148178141
$ws.Range("J9").Value = 0.1986339148178141
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 40.49537033333333
$ws.Range("N9").Value = 121.486111
$ws.Range("O9").Value = 0.4547385869013757
$ws.Range("P9").Value = 0.4547385869013756
$ws.Range("Q9").Value = 210.3045527582729
$ws.Range("R9").Value = 1892.740974824456
$ws.Range("S9").Value = 0.09032650573494103
$ws.Range("T9").Value = 0.090326505734941

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.193298666666666
$ws.Range("H10").Value = 15.579896
$ws.Range("I10").Value = 0.1986339148178141
$ws.Range("J10").Value = 0.1986339148178141
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.218724
$ws.Range("N10").Value = 39.656172
$ws.Range("O10").Value = 0.1484382985738831
$ws.Range("P10").Value = 0.148438298573883
$ws.Range("Q10").Value = 68.64878172423467
$ws.Range("R10").Value = 617.8390355181119
$ws.Range("S10").Value = 0.02948488035462595
$ws.Range("T10").Value = 0.02948488035462594

